$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel;
# mark them as Text before assigning so the literal string is preserved.
$textForceRefs = @("D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D18", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.279.14'
$ws.Range('E2').Value = '  +1.67%  '
$ws.Range('D3').Value = '1.863.48'
$ws.Range('E3').Value = '  +1.46%  '
$ws.Range('E4').Value = '  +1.38%  '
$ws.Range('D5').Value = '312.82'
$ws.Range('E5').Value = '  +1.28%  '
$ws.Range('E6').Value = '  +1.35%  '
$ws.Range('D7').Value = '0.4805'
$ws.Range('E7').Value = '  +1.89%  '
$ws.Range('D8').Value = '0.3732'
$ws.Range('E8').Value = '  +2.16%  '
$ws.Range('D9').Value = '0.07435'
$ws.Range('E9').Value = '  +4.14%  '
$ws.Range('D10').Value = '0.9384'
$ws.Range('E10').Value = '  +2.23%  '
$ws.Range('D11').Value = '20.82'
$ws.Range('E11').Value = '  +6.74%  '
$ws.Range('D12').Value = '0.07876'
$ws.Range('E12').Value = '  +2.95%  '
$ws.Range('D13').Value = '1.868.12'
$ws.Range('E13').Value = '  +4.92%  '
$ws.Range('D14').Value = '5.439'
$ws.Range('E14').Value = '  +3.06%  '
$ws.Range('D15').Value = '6.546'
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('D16').Value = '90.35'
$ws.Range('E16').Value = '  +2.77%  '
$ws.Range('E17').Value = '  +1.35%  '
$ws.Range('D18').Value = '0.000008784'
$ws.Range('E18').Value = '  +1.91%  '
$ws.Range('D20').Value = '14.86'
$ws.Range('E20').Value = '  +2.83%  '
$ws.Range('D21').Value = '27.315.49'
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('D22').Value = '5.120'
$ws.Range('E22').Value = '  +2.34%  '
$ws.Range('D23').Value = '10.70'
$ws.Range('E23').Value = '  +1.06%  '
$ws.Range('D24').Value = '1.958'
$ws.Range('E24').Value = '  +2.02%  '
$ws.Range('D25').Value = '154.16'
$ws.Range('E25').Value = '  +1.64%  '
$ws.Range('D26').Value = '18.56'
$ws.Range('E26').Value = '  +2.05%  '
$ws.Range('D27').Value = '2.010'
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('D28').Value = '116.01'
$ws.Range('E28').Value = '  +1.71%  '
$ws.Range('D29').Value = '5.002'
$ws.Range('E29').Value = '  +2.70%  '
$ws.Range('D30').Value = '0.08913'
$ws.Range('E30').Value = '  +1.13%  '
$ws.Range('D31').Value = '3.344'
$ws.Range('E31').Value = '  +4.30%  '
$ws.Range('D32').Value = '1.200'
$ws.Range('E32').Value = '  +2.40%  '
$ws.Range('D33').Value = '4.579'
$ws.Range('E33').Value = '  +2.57%  '
$ws.Range('D34').Value = '0.7466'
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('D35').Value = '2.681'
$ws.Range('E35').Value = '  -2.24%  '
$ws.Range('D36').Value = '0.02055'
$ws.Range('E36').Value = '  +5.94%  '
$ws.Range('E37').Value = '  +3.44%  '
$ws.Range('D38').Value = '0.05291'
$ws.Range('E38').Value = '  +1.63%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '2.998'
$ws.Range('E39').Value = '  +1.38%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.5379'
$ws.Range('E40').Value = '  +3.90%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '7.116'
$ws.Range('E41').Value = '  +2.36%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = '0.1540'
$ws.Range('E42').Value = '  +2.00%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '8.399'
$ws.Range('E43').Value = '  +3.27%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '10.62'
$ws.Range('E44').Value = '  +1.33%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.4818'
$ws.Range('E45').Value = '  +2.76%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').Value = '1.019'
$ws.Range('E46').Value = '  +1.49%  '
$ws.Range('D47').Value = '1.664'
$ws.Range('E47').Value = '  +4.66%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '103.11'
$ws.Range('E48').Value = '  +1.37%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '66.75'
$ws.Range('E49').Value = '  +3.03%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.06084'
$ws.Range('E50').Value = '  +0.88%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').Value = '0.9014'
$ws.Range('E51').Value = '  +1.97%  '
